{"js": "// Change the bold heading \"REALIZACI\u00d3N DEL PROTOCOLO\" (part of the\n// \"FECHA DE REALIZACI\u00d3N DEL PROTOCOLO:\" label) to\n// \"REALIZACI\u00d3N DEL CONSENTIMIENTO\", keeping the existing bold\n// Book Antiqua / size 20 character formatting intact.\nconst body = context.document.body;\n\nconst results = body.search(\"REALIZACI\u00d3N DEL PROTOCOLO\", { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error('Could not find \"REALIZACI\u00d3N DEL PROTOCOLO\" in the document body.');\n}\n\n// Replace the text in place; formatting (bold / font / size) of the run\n// carries over automatically since insertText with \"Replace\" preserves\n// the existing run formatting.\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(\"REALIZACI\u00d3N DEL CONSENTIMIENTO\", \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# Change the bold heading \"REALIZACI\u00d3N DEL PROTOCOLO\" (part of the\n# \"FECHA DE REALIZACI\u00d3N DEL PROTOCOLO:\" label) to\n# \"REALIZACI\u00d3N DEL CONSENTIMIENTO\", keeping the existing bold\n# Book Antiqua / size 20 character formatting intact.\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"REALIZACI\u00d3N DEL PROTOCOLO\"\n$find.Replacement.Text = \"REALIZACI\u00d3N DEL CONSENTIMIENTO\"\n$find.Forward = $true\n$find.Wrap = 1\n$find.Format = $false\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2)\n"}
